# Hjemme passive updated meanEMG legmaxROM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header/index row) updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updates - B2 and D2 get new values, C2 and E2 are cleared
$ws.Range("B2").Value = 11.467718355161836
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 16.666301382845049
$ws.Range("E2").ClearContents()

# Row 3 (STR) updates
$ws.Range("B3").Value = 10.83486683656362
$ws.Range("C3").Value = -3.105531684919832
$ws.Range("D3").Value = 13.457693562100637
$ws.Range("E3").Value = -6.2146140962329639

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
